# Applies the "additional_parameters" template update:
#  - move "upper_allowable_time_bounds" ahead of "service_history"
#  - insert a new "additional_parameters" sheet (before the hidden
#    "allowable_values" sheet) containing a parameter/value table
#  - leave that new sheet active/selected, matching the authored file

$wb = $excel.ActiveWorkbook

# --- 1. Reorder sheets: upper_allowable_time_bounds moves before service_history ---
$wsServiceHistory = $wb.Worksheets.Item("service_history")
$wsUpperBounds = $wb.Worksheets.Item("upper_allowable_time_bounds")
$wsUpperBounds.Move($wsServiceHistory)

# --- 2. Insert the new worksheet right before the hidden allowable_values sheet ---
$wsAllowableValues = $wb.Worksheets.Item("allowable_values")
$wsParams = $wb.Worksheets.Add($wsAllowableValues)
$wsParams.Name = "additional_parameters"

# --- 3. Column widths (bestFit, matches authored widths) ---
$wsParams.Columns.Item(1).ColumnWidth = 34.85546875
$wsParams.Columns.Item(2).ColumnWidth = 10.42578125

# --- 4. Header row ---
$wsParams.Range("A1").Value = "parameter"
$wsParams.Range("B1").Value = "value"

# --- 5. Parameter rows ---
$wsParams.Range("A2").Value = "number_of_runs"
$wsParams.Range("B2").Value = 12

$wsParams.Range("A3").Value = "simulation_duration_days"
$wsParams.Range("B3").Value = 730

$wsParams.Range("A4").Value = "simulation_warm_up_duration_hours"
$wsParams.Range("B4").Value = 0

$wsParams.Range("A5").Value = "simulation_start_date"
$wsParams.Range("B5").Value = 44927
$wsParams.Range("B5").NumberFormat = "mm-dd-yy"

$wsParams.Range("A6").Value = "simulation_start_time"
$wsParams.Range("B6").Value = 0.33333333333333331
$wsParams.Range("B6").NumberFormat = "h:mm"

$wsParams.Range("A7").Value = "master_random_seed"
$wsParams.Range("B7").Value = 42

$wsParams.Range("A8").Value = "activity_duration_multiplier"
$wsParams.Range("B8").Value = 1

# --- 6. Turn the range into a table, matching the workbook's other tables ---
$tbl = $wsParams.ListObjects.Add(1, $wsParams.Range("A1:B21"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight8"

# --- 7. Make the new sheet the active/selected one, like the authored file ---
$wsParams.Activate()
$wsParams.Range("A1").Select()
